# Generate Report for Handback
# Populates handback status/metadata for the two localized files (zh-cn, de-de)
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Update the "Status" columns so every row reflects that the handback is
#    complete and in sync with en-US.
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (I) / "Latest Handback File"
#    (J) with the target markdown file and the handback xliff for each row,
#    and add the corresponding hyperlink to the target file (mirrors column A).
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/178b2a701da26ba21b2d3751489c2995d07045b7/e2e/12040c64-1c18-422e-ab12-e661cbc401b8.md",
    "",
    "",
    "12040c64-1c18-422e-ab12-e661cbc401b8.md"
) | Out-Null
$wsZhCn.Range("J2").Value = "12040c64-1c18-422e-ab12-e661cbc401b8.4938b1e8b6d73aa270f8644325b904a853757f94.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/178b2a701da26ba21b2d3751489c2995d07045b7/e2e/5d3e1a15-ee8c-4620-8cdf-9ca9e09c2bb7.md",
    "",
    "",
    "5d3e1a15-ee8c-4620-8cdf-9ca9e09c2bb7.md"
) | Out-Null
$wsZhCn.Range("J3").Value = "5d3e1a15-ee8c-4620-8cdf-9ca9e09c2bb7.4eb7f3a252dbfcbd45ad2781d9fc93480653a9c1.zh-cn.xlf"

# Latest Handback DateTime for zh-cn
$wsZhCn.Range("K2").Value = "2016-08-16 12:46:45"
$wsZhCn.Range("K3").Value = "2016-08-16 12:46:45"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same as above, plus its own (later) handback timestamp.
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/178b2a701da26ba21b2d3751489c2995d07045b7/e2e/12040c64-1c18-422e-ab12-e661cbc401b8.md",
    "",
    "",
    "12040c64-1c18-422e-ab12-e661cbc401b8.md"
) | Out-Null
$wsDeDe.Range("J2").Value = "12040c64-1c18-422e-ab12-e661cbc401b8.4938b1e8b6d73aa270f8644325b904a853757f94.de-de.xlf"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/178b2a701da26ba21b2d3751489c2995d07045b7/e2e/5d3e1a15-ee8c-4620-8cdf-9ca9e09c2bb7.md",
    "",
    "",
    "5d3e1a15-ee8c-4620-8cdf-9ca9e09c2bb7.md"
) | Out-Null
$wsDeDe.Range("J3").Value = "5d3e1a15-ee8c-4620-8cdf-9ca9e09c2bb7.4eb7f3a252dbfcbd45ad2781d9fc93480653a9c1.de-de.xlf"

# Latest Handback DateTime for de-de
$wsDeDe.Range("K2").Value = "2016-08-16 12:46:53"
$wsDeDe.Range("K3").Value = "2016-08-16 12:46:53"

# ---------------------------------------------------------------------------
# 4. Widen the columns that now hold longer text (Status columns on Overview,
#    Status column on the language sheets, and the newly populated Latest
#    Target File / Latest Handback File columns).
# ---------------------------------------------------------------------------
$wideStatusWidth = 29.166666666666664   # renders as OOXML width ~30 (closest reachable to 29.9777)
$wideFileWidth   = 39.2                 # renders as OOXML width 40

# Overview: columns E (5) and F (6)
$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth

# zh-cn: column C (3), I (9), J (10)
$wsZhCn.Columns.Item(3).ColumnWidth  = $wideStatusWidth
$wsZhCn.Columns.Item(9).ColumnWidth  = $wideFileWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideFileWidth

# de-de: column C (3), I (9), J (10)
$wsDeDe.Columns.Item(3).ColumnWidth  = $wideStatusWidth
$wsDeDe.Columns.Item(9).ColumnWidth  = $wideFileWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideFileWidth
